# Add a new book entry to the wishlist (row 10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "El reino del dragón de oro"
$ws.Range("B10").Value = "Isabel Allende"

# Mirror the empty Editorial cell style/shape used by the rows above
# (C8/C9 are present but empty) so C10 matches that same pattern.
$ws.Range("C9").Copy($ws.Range("C10"))
